# Adds a new "Mediocre medio" match result (Alba-Luis beats Alberto-Esperanza,
# 6-4,6-1 / 4-6,1-6) to the resultados sheet, rolls the historial_partidos
# dates forward one day, appends the two corresponding historial_partidos
# rows, and updates the derived classification stats on clasificacion_auto
# and clasificacion for the two pairs involved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) resultados: new row 24 with the new match result
# ---------------------------------------------------------------------
$resultados = $wb.Worksheets.Item("resultados")
$resultados.Cells.Item(24, 1).Value = "Mediocre medio"
$resultados.Cells.Item(24, 2).Value = "1ª vuelta"
$resultados.Cells.Item(24, 3).Value = "Alba-Luis"
$resultados.Cells.Item(24, 4).Value = "Alberto-Esperanza"
$resultados.Cells.Item(24, 5).Value = "6-4,6-1"
$resultados.Cells.Item(24, 6).Value = "4-6,1-6"
$resultados.Range("G24").Select()

# ---------------------------------------------------------------------
# 2) clasificacion_auto: update stats for Alberto-Esperanza (row 17) and
#    Alba-Luis (row 18)
# ---------------------------------------------------------------------
$clasifAuto = $wb.Worksheets.Item("clasificacion_auto")

$clasifAuto.Cells.Item(17, 5).Value = 2
$clasifAuto.Cells.Item(17, 8).Value = 1
$clasifAuto.Cells.Item(17, 10).Value = 2
$clasifAuto.Cells.Item(17, 11).Value = 16
$clasifAuto.Cells.Item(17, 12).Value = 17

$clasifAuto.Cells.Item(18, 4).Value = 3
$clasifAuto.Cells.Item(18, 5).Value = 4
$clasifAuto.Cells.Item(18, 6).Value = 1
$clasifAuto.Cells.Item(18, 8).Value = 3
$clasifAuto.Cells.Item(18, 9).Value = 2
$clasifAuto.Cells.Item(18, 10).Value = 6
$clasifAuto.Cells.Item(18, 11).Value = 20
$clasifAuto.Cells.Item(18, 12).Value = 41

# ---------------------------------------------------------------------
# 3) clasificacion: update stats for Alberto-Esperanza (row 10) and
#    Alba-Luis (row 11)
# ---------------------------------------------------------------------
$clasif = $wb.Worksheets.Item("clasificacion")

$clasif.Cells.Item(10, 5).Value = 2
$clasif.Cells.Item(10, 8).Value = 1
$clasif.Cells.Item(10, 10).Value = 2

$clasif.Cells.Item(11, 4).Value = 3
$clasif.Cells.Item(11, 5).Value = 4
$clasif.Cells.Item(11, 6).Value = 1
$clasif.Cells.Item(11, 8).Value = 3
$clasif.Cells.Item(11, 9).Value = 2
$clasif.Cells.Item(11, 10).Value = 6

# ---------------------------------------------------------------------
# 4) historial_partidos: bump every existing date (rows 2-45) forward one
#    day, then append the two new match rows (46, 47)
# ---------------------------------------------------------------------
$historial = $wb.Worksheets.Item("historial_partidos")

for ($r = 2; $r -le 45; $r++) {
    $historial.Cells.Item($r, 1).Value = 45975
}

$historial.Cells.Item(46, 1).Value = 45975
$historial.Cells.Item(46, 2).Value = "mediocre medio"
$historial.Cells.Item(46, 3).Value = "1ª vuelta"
$historial.Cells.Item(46, 4).Value = "Alba-Luis"
$historial.Cells.Item(46, 5).Value = "Gana"
$historial.Cells.Item(46, 6).Value = 2
$historial.Cells.Item(46, 7).Value = 0
$historial.Cells.Item(46, 8).Value = 3
$historial.Cells.Item(46, 9).Value = 2
$historial.Cells.Item(46, 10).Value = 3
$historial.Cells.Item(46, 11).Value = 1
$historial.Cells.Item(46, 12).Value = 0
$historial.Cells.Item(46, 13).Value = 1

$historial.Cells.Item(47, 1).Value = 45975
$historial.Cells.Item(47, 2).Value = "mediocre medio"
$historial.Cells.Item(47, 3).Value = "1ª vuelta"
$historial.Cells.Item(47, 4).Value = "Alberto-Esperanza"
$historial.Cells.Item(47, 5).Value = "Pierde"
$historial.Cells.Item(47, 6).Value = 0
$historial.Cells.Item(47, 7).Value = 2
$historial.Cells.Item(47, 8).Value = 0
$historial.Cells.Item(47, 9).Value = 4
$historial.Cells.Item(47, 10).Value = 3
$historial.Cells.Item(47, 11).Value = 1
$historial.Cells.Item(47, 12).Value = 0
$historial.Cells.Item(47, 13).Value = 3
